# Workbook was deployed under a new path; refresh the sheet name that was
# left over from the old template (sheet1 -> Sheet1) and move the saved
# selection to where the user was last working (C21) after the deploy.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1"

$ws.Activate()
$ws.Range("C21").Select()
